# Update "想去人数" (interest count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 1322
$wsExhibit.Range("F3").Value  = 1205
$wsExhibit.Range("F4").Value  = 14513
$wsExhibit.Range("F5").Value  = 17291
$wsExhibit.Range("F7").Value  = 144
$wsExhibit.Range("F8").Value  = 56
$wsExhibit.Range("F10").Value = 206
$wsExhibit.Range("F16").Value = 40
$wsExhibit.Range("F17").Value = 20
$wsExhibit.Range("F18").Value = 129
$wsExhibit.Range("F23").Value = 61
$wsExhibit.Range("F25").Value = 7118
$wsExhibit.Range("F30").Value = 5831
$wsExhibit.Range("F31").Value = 51
$wsExhibit.Range("F35").Value = 218
$wsExhibit.Range("F36").Value = 5011

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 1322
$wsAll.Range("F3").Value  = 1205
$wsAll.Range("F4").Value  = 14513
$wsAll.Range("F5").Value  = 17291
$wsAll.Range("F7").Value  = 144
$wsAll.Range("F8").Value  = 56
$wsAll.Range("F10").Value = 206
$wsAll.Range("F16").Value = 40
$wsAll.Range("F17").Value = 20
$wsAll.Range("F18").Value = 129
$wsAll.Range("F24").Value = 61
$wsAll.Range("F26").Value = 7118
$wsAll.Range("F32").Value = 5831
$wsAll.Range("F33").Value = 51
$wsAll.Range("F37").Value = 218
$wsAll.Range("F38").Value = 5011
